# Apply Optuna-attempt forecast value updates to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 234
$ws1.Range("H2").Value = 6.15
$ws1.Range("L2").Value = 0.8100000000000001

# Row 3
$ws1.Range("D3").Value = 249
$ws1.Range("H3").Value = 4.84
$ws1.Range("L3").Value = 1.17

# Row 4
$ws1.Range("D4").Value = 243
$ws1.Range("H4").Value = 3.94
$ws1.Range("L4").Value = 1.11

# Row 5
$ws1.Range("H5").Value = 2.8

# Row 6
$ws1.Range("H6").Value = 1.78
$ws1.Range("J6").Value = "Normal"
$ws1.Range("L6").Value = 1.07

# Row 7
$ws1.Range("H7").Value = 0.84
$ws1.Range("I7").Value = "Low"
$ws1.Range("L7").Value = 0.85

# Row 8
$ws1.Range("L8").Value = 0.8

# Row 9
$ws1.Range("L9").Value = 1.11

# Row 10
$ws1.Range("L10").Value = 1

# Row 11
$ws1.Range("L11").Value = 1.19

# Row 12
$ws1.Range("L12").Value = 0.9

# Row 13
$ws1.Range("L13").Value = 0.88

# Row 14
$ws1.Range("L14").Value = 0.99

# Row 15
$ws1.Range("L15").Value = 0.96

# Row 16
$ws1.Range("L16").Value = 0.99

# Row 17
$ws1.Range("L17").Value = 1.11

# --- Sheet 2: "Summary" ---
# These cells hold numeric-looking text (stored as strings in the
# workbook), so force Text format before assigning so COM doesn't
# silently coerce the value to a Number.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "3759"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "1968"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "981"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "259"
